# feat: add 2022-Q4 data
#
# - Insert a new worksheet "2022-Q4" right after "总计" and before "2022-Q3",
#   populated with the per-fund holding detail table for that quarter.
# - Insert a new row into the "总计" (totals) sheet summarizing the new
#   quarter, shifting the existing rows down.

$wb = $excel.ActiveWorkbook

function Set-TextCell($range, [string]$val) {
    # Force text storage even for numeric-looking strings (e.g. "40.90",
    # "001411") so fund codes keep leading zeros and decimals keep their
    # trailing zeros, then drop the residual number-format style so the
    # cell is left with the default (unstyled) look - matching how the
    # rest of the workbook stores this text data.
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.ClearFormats()
}

function Set-NumCell($range, $val) {
    # Plain numeric cell, no special styling.
    $range.ClearFormats()
    $range.Value = $val
}

function Set-HeaderCell($range, [string]$val) {
    # Bold, bordered, centered header cell (text).
    $range.ClearFormats()
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1
}

function Set-IndexCell($range, $val) {
    # Column-A row-index cells: bold + bordered + centered, numeric value.
    $range.ClearFormats()
    $range.Value = $val
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet before the existing "2022-Q3" sheet.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Header row
Set-HeaderCell $q4.Range("B1") "基金代码"
Set-HeaderCell $q4.Range("C1") "基金名称"
Set-HeaderCell $q4.Range("D1") "基金规模"
Set-HeaderCell $q4.Range("E1") "股票总仓位"
Set-HeaderCell $q4.Range("F1") "仓位占比"
Set-HeaderCell $q4.Range("G1") "持有市值(亿元)"
Set-HeaderCell $q4.Range("H1") "仓位排名"

# Fund detail rows
$rows = @(
    @{ idx=0; code="320003"; name="诺安先锋混合A";               size="40.90"; pos="76.87"; pct="2.75"; mv="1.1248"; rank=8 },
    @{ idx=1; code="001411"; name="诺安创新驱动灵活配置混合A";   size="7.22";  pos="91.15"; pct="4.56"; mv="0.3292"; rank=4 },
    @{ idx=2; code="002051"; name="诺安创新驱动灵活配置混合C";   size="4.60";  pos="91.15"; pct="4.56"; mv="0.2098"; rank=4 },
    @{ idx=3; code="001706"; name="诺安积极回报灵活配置混合A";   size="0.52";  pos="93.31"; pct="7.38"; mv="0.0384"; rank=8 },
    @{ idx=4; code="012621"; name="诺安先锋混合C";               size="1.09";  pos="76.87"; pct="2.75"; mv="0.0300"; rank=8 },
    @{ idx=5; code="012847"; name="诺安积极回报灵活配置混合C";   size="0.18";  pos="93.31"; pct="7.38"; mv="0.0133"; rank=8 }
)

$r = 2
foreach ($row in $rows) {
    Set-IndexCell $q4.Range("A$r") $row.idx
    Set-TextCell  $q4.Range("B$r") $row.code
    Set-TextCell  $q4.Range("C$r") $row.name
    Set-TextCell  $q4.Range("D$r") $row.size
    Set-TextCell  $q4.Range("E$r") $row.pos
    Set-TextCell  $q4.Range("F$r") $row.pct
    Set-TextCell  $q4.Range("G$r") $row.mv
    Set-NumCell   $q4.Range("H$r") $row.rank
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Insert the new summary row into "总计", shifting existing data down.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Copy the (now-shifted) row 3's formatting back onto the fresh row 2 so
# the new row matches the sheet's existing look exactly (bold/bordered
# index cell in column A, unstyled data cells in B:D) instead of whatever
# Excel's row-insert auto-formatting guessed.
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)   # xlPasteFormats

$total.Range("A2").Value = 0
Set-TextCell $total.Range("B2") "2022-Q4"
$total.Range("C2").Value = 6
$total.Range("D2").Value = 1.75
